$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current last data row (row 65), shifting the
# existing row 65 data down to row 66, and populate the new row with the
# latest weekly price entry.
$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65, 1).Value2 = 1
$ws.Cells.Item(65, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(65, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(65, 4).Value2 = 45008
$ws.Cells.Item(65, 5).Value2 = 15
$ws.Cells.Item(65, 6).Value2 = 100112028
$ws.Cells.Item(65, 7).Value2 = "Sandia"
$ws.Cells.Item(65, 8).Value2 = "Sin especificar"
$ws.Cells.Item(65, 9).Value2 = "Segunda"
$ws.Cells.Item(65, 10).Value2 = 400
$ws.Cells.Item(65, 11).Value2 = 450
$ws.Cells.Item(65, 12).Value2 = 480
$ws.Cells.Item(65, 13).Value2 = 465
$ws.Cells.Item(65, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(65, 15).Value2 = "Perú"
$ws.Cells.Item(65, 16).Value2 = 465
$ws.Cells.Item(65, 17).Value2 = 1
$ws.Cells.Item(65, 18).Value2 = "Hortaliza"
